$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Septiembre de 2020 a las 16:38"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 6879675
$ws.Range("C4").Value = 5079
$ws.Range("D4").Value = 4156240
$ws.Range("E4").Value = 2521131
$ws.Range("G4").Value = 91
$ws.Range("H4").Value = 202304

# Row 5: 'India' -> 'India'
$ws.Range("B5").Value = 5228478
$ws.Range("C5").Value = 15792
$ws.Range("D5").Value = 4125742
$ws.Range("E5").Value = 1018231
$ws.Range("G5").Value = 101
$ws.Range("H5").Value = 84505

# Row 13: 'Argentina' -> 'Argentina'
$ws.Range("D13").Value = 467286
$ws.Range("E13").Value = 121936
$ws.Range("G13").Value = 31
$ws.Range("H13").Value = 12491

# Row 14: 'Chile' -> 'Chile'
$ws.Range("B14").Value = 442827
$ws.Range("C14").Value = 1677
$ws.Range("D14").Value = 415981
$ws.Range("E14").Value = 14647
$ws.Range("G14").Value = 57
$ws.Range("H14").Value = 12199

# Row 25: 'Alemania' -> 'Alemania'
$ws.Range("B25").Value = 270123
$ws.Range("C25").Value = 1081
$ws.Range("E25").Value = 19361

# Row 27: 'Israel' -> 'Israel'
$ws.Range("B27").Value = 179071
$ws.Range("C27").Value = 3815
$ws.Range("D27").Value = 130024
$ws.Range("E27").Value = 47851
$ws.Range("G27").Value = 27
$ws.Range("H27").Value = 1196

# Row 51: 'Etiopia' -> 'Portugal'
$ws.Range("A51").Value = "Portugal"
$ws.Range("B51").Value = 67176
$ws.Range("C51").Value = 780
$ws.Range("D51").Value = 45053
$ws.Range("E51").Value = 20229
$ws.Range("G51").Value = 6
$ws.Range("H51").Value = 1894

# Row 52: 'Portugal' -> 'Etiopia'
$ws.Range("A52").Value = "Etiopia"
$ws.Range("B52").Value = 66913
$ws.Range("D52").Value = 27085
$ws.Range("E52").Value = 38768
$ws.Range("H52").Value = 1060

# Row 59: 'Uzbekistan' -> 'Uzbekistan'
$ws.Range("B59").Value = 50253
$ws.Range("C59").Value = 626
$ws.Range("D59").Value = 46527
$ws.Range("E59").Value = 3307
$ws.Range("G59").Value = 6
$ws.Range("H59").Value = 419

# Row 69: 'Austria' -> 'Kenia'
$ws.Range("A69").Value = "Kenia"
$ws.Range("B69").Value = 36724
$ws.Range("C69").Value = 148
$ws.Range("D69").Value = 23709
$ws.Range("E69").Value = 12369
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 646

# Row 70: 'Kenia' -> 'Austria'
$ws.Range("A70").Value = "Austria"
$ws.Range("B70").Value = 36661
$ws.Range("C70").Value = 808
$ws.Range("D70").Value = 28451
$ws.Range("E70").Value = 7447
$ws.Range("G70").Value = 5
$ws.Range("H70").Value = 763

# Row 72: 'Serbia' -> 'Serbia'
$ws.Range("D72").Value = 31411
$ws.Range("E72").Value = 607

# Row 125: 'Jamaica' -> 'Birmania'
$ws.Range("A125").Value = "Birmania"
$ws.Range("B125").Value = 4467
$ws.Range("C125").Value = 424
$ws.Range("D125").Value = 1130
$ws.Range("E125").Value = 3267
$ws.Range("G125").Value = 10
$ws.Range("H125").Value = 70

# Row 126: 'Birmania' -> 'Jamaica'
$ws.Range("A126").Value = "Jamaica"
$ws.Range("B126").Value = 4374
$ws.Range("C126").Value = 0
$ws.Range("D126").Value = 1225
$ws.Range("E126").Value = 3098
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 51

# Row 201: 'Guam' -> 'Bonaire, San Eustaquio y Saba'
$ws.Range("A201").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B201").Value = 36
$ws.Range("C201").Value = 4
$ws.Range("D201").Value = 17
$ws.Range("E201").Value = 18

# Row 202: 'Bonaire, San Eustaquio y Saba' -> 'Guam'
$ws.Range("A202").Value = "Guam"
$ws.Range("D202").Value = 0
$ws.Range("E202").Value = 31

# Row 214: 'Islas Malvinas' -> 'Montserrat'
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Row 215: 'Montserrat' -> 'Islas Malvinas'
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
